$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Round the coordinate values in Q2 and R2 to whole numbers.
$ws.Range("Q2").Value = 369129
$ws.Range("R2").Value = 6363255

# Remove the time values (Starttid / Sluttid) entirely.
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
